$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# vehicleInfo (sheet2): replace the placeholder "CX-7 / SUV (4CYL 4x2)" row
# with a "null" type marker on row 2 and new 2012 BMW 328 / 2014 Kia Forte
# vehicle rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("vehicleInfo")
$ws2.Range("D2").Value = "null"

$ws2.Range("A3").Value = "'2012"
$ws2.Range("B3").Value = "BMW"
$ws2.Range("C3").Value = "'328"
$ws2.Range("D3").Value = "2DR 6CYL"

$ws2.Range("A4").Value = "'2014"
$ws2.Range("B4").Value = "Kia"
$ws2.Range("C4").Value = "Forte"
$ws2.Range("D4").Value = "2DR 4CYL"
$ws2.Range("E4").Value = "Personal (to/from work or school, errands, pleasure)"
$ws2.Range("F4").Value = "Finance"
$ws2.Range("G4").Value = "1 year - 3 years"

$ws2.Range("D32").Select()

# ---------------------------------------------------------------------------
# removeVehicle (new sheet3)
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "removeVehicle"
$ws3.Range("A1").Value = "make"
$ws3.Range("A2").Value = "BMW"
$ws3.Range("A2").Select()

# ---------------------------------------------------------------------------
# driverInformation (new sheet4)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "driverInformation"
$ws4.Range("A1").Value = "gender"
$ws4.Range("B1").Value = "selectMaritalStatus"
$ws4.Range("C1").Value = "selectEducation"
$ws4.Range("D1").Value = "selectEmployment"
$ws4.Range("F1").Value = "socialSecurityNumber"
$ws4.Range("G1").Value = "primaryResidency"
$ws4.Range("H1").Value = "hasPriorAddress"
$ws4.Range("I1").Value = "licensed"
$ws4.Range("J1").Value = "accident"
$ws4.Range("K1").Value = "ticket"

$ws4.Range("A2").Value = "Female"
$ws4.Range("B2").Value = "Single"
$ws4.Range("C2").Value = "College degree"
$ws4.Range("D2").Value = "Employed"

$ws4.Range("E1").Value = "enterOccupation"
$ws4.Range("E2").Value = "Quality Assurance Tester"

$ws4.Range("F2").Value = 123456789
$ws4.Range("G2").Value = "Own home"
$ws4.Range("H2").Value = "No"
$ws4.Range("I2").Value = "At least 2 years, but less than 3 years"
$ws4.Range("J2").Value = "No"
$ws4.Range("K2").Value = "No"

$ws4.Range("A2").Select()

# ---------------------------------------------------------------------------
# AdditionalDetail (new sheet5)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "AdditionalDetail"
$ws5.Range("A1").Value = "insuranceToday"
$ws5.Range("B1").Value = "insuranceLastMonth"
$ws5.Range("C1").Value = "nonAutoPolicy"
$ws5.Range("D1").Value = "emailAddress"
$ws5.Range("E1").Value = "totalResident"
$ws5.Range("F1").Value = "residence"

$ws5.Range("A2").Value = "No"
$ws5.Range("B2").Value = "No"
$ws5.Range("C2").Value = "No"

$ws5.Range("D2").Value = "stha@gmail.com"
$ws5.Hyperlinks.Add($ws5.Range("D2"), "mailto:stha@gmail.com")

$ws5.Range("E2").Value = 1
$ws5.Range("F2").Value = "Less than 1 year"

$ws5.Range("A2").Select()
